$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update semantic metadata for estado-civil: it is now curated as a measure, not a dimension.
$ws.Range("A2").Value = "iaest-measure:estado-civil"

# The "aragon" column no longer has its own dimension concept; it now reuses the
# standard SDMX reference-area dimension concept.
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# estado-civil's row-3 metadata type switches from "dim" to "medida" (measure).
$ws.Range("A3").Value = "medida"

# estado-civil's data type switches from skos:Concept to xsd:int (now a measure value).
$ws.Range("A4").Value = "xsd:int"

# The "aragon" column's URI template changes from skos:Concept to the new URI-Comunidad template.
$ws.Range("E4").Value = "URI-Comunidad"

# Remove the obsolete mapping file references row entirely (no longer needed with curated dimensions).
$ws.Rows.Item(5).Delete()
